$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking / percentage text columns to Text format so the
# values round-trip as strings (matching the inlineStr source data),
# rather than being auto-converted to numbers/percentages.
$textRanges = $ws.Range("D2:D51,E2:E51,G2:G51")
foreach ($area in $textRanges.Areas) {
    $area.NumberFormat = "@"
}

# Apply the updated cell values from the diff.
$ws.Range("D2").Value = "332.18"
$ws.Range("E2").Value = "1.35%"
$ws.Range("G2").Value = "6"
$ws.Range("D3").Value = "45.89"
$ws.Range("E3").Value = "4.63%"
$ws.Range("G3").Value = "6"
$ws.Range("D4").Value = "5.640"
$ws.Range("E4").Value = "2.50%"
$ws.Range("G4").Value = "6"
$ws.Range("D5").Value = "0.08374"
$ws.Range("E5").Value = "4.51%"
$ws.Range("G5").Value = "6"
$ws.Range("D6").Value = "2.041"
$ws.Range("E6").Value = "1.53%"
$ws.Range("G6").Value = "6"
$ws.Range("D7").Value = "4.484"
$ws.Range("E7").Value = "3.79%"
$ws.Range("G7").Value = "6"
$ws.Range("D8").Value = "0.9940"
$ws.Range("E8").Value = "4.75%"
$ws.Range("G8").Value = "6"
$ws.Range("D9").Value = "2.582"
$ws.Range("E9").Value = "0.50%"
$ws.Range("G9").Value = "6"
$ws.Range("D10").Value = "0.1155"
$ws.Range("E10").Value = "2.87%"
$ws.Range("G10").Value = "6"
$ws.Range("D11").Value = "0.1937"
$ws.Range("E11").Value = "4.26%"
$ws.Range("G11").Value = "6"
$ws.Range("D12").Value = "10.42"
$ws.Range("E12").Value = "-1.86%"
$ws.Range("G12").Value = "6"
$ws.Range("D13").Value = "0.09946"
$ws.Range("E13").Value = "0.26%"
$ws.Range("G13").Value = "6"
$ws.Range("D14").Value = "0.04683"
$ws.Range("E14").Value = "1.34%"
$ws.Range("G14").Value = "6"
$ws.Range("D15").Value = "0.1060"
$ws.Range("E15").Value = "-0.60%"
$ws.Range("G15").Value = "6"
$ws.Range("D16").Value = "0.001275"
$ws.Range("E16").Value = "0.38%"
$ws.Range("G16").Value = "6"
$ws.Range("D17").Value = "0.006084"
$ws.Range("E17").Value = "2.42%"
$ws.Range("G17").Value = "6"
$ws.Range("D18").Value = "3.377"
$ws.Range("E18").Value = "0.60%"
$ws.Range("G18").Value = "6"
$ws.Range("D19").Value = "0.3367"
$ws.Range("G19").Value = "6"
$ws.Range("D20").Value = "0.1404"
$ws.Range("E20").Value = "-0.17%"
$ws.Range("G20").Value = "6"
$ws.Range("E21").Value = "4.32%"
$ws.Range("G21").Value = "6"
$ws.Range("D22").Value = "0.04220"
$ws.Range("E22").Value = "3.54%"
$ws.Range("G22").Value = "6"
$ws.Range("D23").Value = "0.001313"
$ws.Range("E23").Value = "4.38%"
$ws.Range("G23").Value = "6"
$ws.Range("D24").Value = "0.004668"
$ws.Range("E24").Value = "7.66%"
$ws.Range("G24").Value = "6"
$ws.Range("D25").Value = "0.0001285"
$ws.Range("E25").Value = "10.86%"
$ws.Range("G25").Value = "6"
$ws.Range("D26").Value = "0.0003753"
$ws.Range("E26").Value = "0.33%"
$ws.Range("G26").Value = "6"
$ws.Range("G27").Value = "6"
$ws.Range("G28").Value = "6"
$ws.Range("G29").Value = "6"
$ws.Range("G30").Value = "6"
$ws.Range("G31").Value = "6"
$ws.Range("G32").Value = "6"
$ws.Range("G33").Value = "6"
$ws.Range("G34").Value = "6"
$ws.Range("G35").Value = "6"
$ws.Range("G36").Value = "6"
$ws.Range("G37").Value = "6"
$ws.Range("D38").Value = "0.02792"
$ws.Range("E38").Value = "8.05%"
$ws.Range("G38").Value = "6"
$ws.Range("E39").Value = "2.21%"
$ws.Range("G39").Value = "6"
$ws.Range("D40").Value = "0.007803"
$ws.Range("E40").Value = "3.57%"
$ws.Range("G40").Value = "6"
$ws.Range("D41").Value = "0.1438"
$ws.Range("E41").Value = "2.97%"
$ws.Range("G41").Value = "6"
$ws.Range("D42").Value = "0.007274"
$ws.Range("E42").Value = "-3.22%"
$ws.Range("G42").Value = "6"
$ws.Range("D43").Value = "0.002018"
$ws.Range("E43").Value = "0.21%"
$ws.Range("G43").Value = "6"
$ws.Range("D44").Value = "0.009050"
$ws.Range("E44").Value = "7.98%"
$ws.Range("G44").Value = "6"
$ws.Range("D45").Value = "0.3408"
$ws.Range("G45").Value = "6"
$ws.Range("D46").Value = "0.00007333"
$ws.Range("E46").Value = "3.22%"
$ws.Range("G46").Value = "6"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").Value = "0.40%"
$ws.Range("G47").Value = "6"
$ws.Range("D48").Value = "0.0005821"
$ws.Range("E48").Value = "0.16%"
$ws.Range("G48").Value = "6"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "0.003504"
$ws.Range("E49").Value = "1.45%"
$ws.Range("G49").Value = "6"
$ws.Range("B50").Value = "CoinbaseStockToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D50").Value = "0.003510"
$ws.Range("E50").Value = "-0.53%"
$ws.Range("G50").Value = "6"
$ws.Range("D51").Value = "0.00002107"
$ws.Range("E51").Value = "0.40%"
$ws.Range("G51").Value = "6"

# Strip the temporary text-format override again so the cells keep their
# original (default) style, matching the source workbook formatting.
foreach ($area in $textRanges.Areas) {
    $area.ClearFormats()
}
